$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistroEntradas")
$ws.Range("B5").Value = "Data do Caixa Realizado"
$rng = $ws.Range("B5:H9")
$lo = $ws.ListObjects.Add(1, $rng, 0, 1)
$lo.Name = "TbRegistroEntradas"
$lo.TableStyle = "TableStyleLight15"
Write-Output $lo.TableStyle
